$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    # Force the cell to remain plain text even when the new value
    # looks numeric (e.g. "1.00"), then restore the default (unstyled)
    # cell style so no formatting is introduced.
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $value
    $ws.Range($cellRef).Style = "Normal"
}

$ws.Range("D2").Value = '61.318.61'
$ws.Range("E2").Value = '  -1.84%  '
$ws.Range("D3").Value = '2.981.82'
$ws.Range("E3").Value = '  -1.34%  '
Set-TextValue "D4" '1.00'
$ws.Range("E4").Value = '  -0.03%  '
Set-TextValue "D5" '595.29'
$ws.Range("E5").Value = '  +1.61%  '
Set-TextValue "D6" '143.93'
$ws.Range("E6").Value = '  -2.58%  '
$ws.Range("D8").Value = '2.982.18'
$ws.Range("E8").Value = '  -1.32%  '
Set-TextValue "D9" '0.513'
$ws.Range("E9").Value = '  -2.67%  '
$ws.Range("E10").Value = '  -1.11%  '
Set-TextValue "D11" '6.05'
$ws.Range("E11").Value = '  +3.24%  '
Set-TextValue "D12" '0.452'
$ws.Range("E12").Value = '  +0.35%  '
$ws.Range("E13").Value = '  -1.39%  '
Set-TextValue "D14" '34.05'
$ws.Range("E14").Value = '  -1.93%  '
$ws.Range("E15").Value = '  +3.16%  '
$ws.Range("D16").Value = '3.470.95'
$ws.Range("E16").Value = '  -1.40%  '
Set-TextValue "D17" '6.91'
$ws.Range("E17").Value = '  -2.49%  '
$ws.Range("D18").Value = '61.276.71'
$ws.Range("E18").Value = '  -1.83%  '
$ws.Range("D19").Value = '2.978.19'
$ws.Range("E19").Value = '  -1.38%  '
Set-TextValue "D20" '444.91'
$ws.Range("E20").Value = '  -3.96%  '
Set-TextValue "D21" '13.90'
$ws.Range("E21").Value = '  -0.74%  '
Set-TextValue "D22" '0.681'
$ws.Range("E22").Value = '  -0.91%  '
$ws.Range("E23").Value = '  -2.13%  '
Set-TextValue "D24" '80.83'
$ws.Range("E24").Value = '  -1.02%  '
Set-TextValue "D25" '10.74'
$ws.Range("E25").Value = '  +5.17%  '
$ws.Range("E26").Value = '  -3.92%  '
Set-TextValue "D27" '12.02'
$ws.Range("E27").Value = '  -2.77%  '
$ws.Range("E28").Value = '  +0.13%  '
$ws.Range("E29").Value = '  +2.36%  '
$ws.Range("E30").Value = '  -0.01%  '
Set-TextValue "D31" '7.19'
$ws.Range("E31").Value = '  +1.01%  '
Set-TextValue "D32" '2.05'
$ws.Range("E32").Value = '  -2.94%  '
Set-TextValue "D33" '27.14'
$ws.Range("E33").Value = '  -4.84%  '
$ws.Range("E34").Value = '  +0.15%  '
$ws.Range("D35").Value = '0.0₃0805'
$ws.Range("E35").Value = '  -0.08%  '
$ws.Range("E36").Value = '  -1.71%  '
Set-TextValue "D37" '5.75'
$ws.Range("E37").Value = '  -0.71%  '
Set-TextValue "D38" '50.19'
$ws.Range("E38").Value = '  -0.48%  '
Set-TextValue "D39" '8.96'
$ws.Range("E39").Value = '  -1.70%  '
Set-TextValue "D40" '2.01'
$ws.Range("E40").Value = '  -5.32%  '
Set-TextValue "D41" '0.124'
$ws.Range("E41").Value = '  +9.14%  '
Set-TextValue "D42" '2.82'
$ws.Range("E42").Value = '  -4.20%  '
Set-TextValue "D43" '386.74'
$ws.Range("E43").Value = '  -1.64%  '
Set-TextValue "D44" '39.41'
$ws.Range("E44").Value = '  +5.98%  '
$ws.Range("E45").Value = '  -3.18%  '
$ws.Range("E46").Value = '  -2.90%  '
$ws.Range("D47").Value = '2.682.03'
$ws.Range("E47").Value = '  -2.76%  '
Set-TextValue "D48" '131.19'
$ws.Range("E48").Value = '  +1.98%  '
Set-TextValue "D50" '0.107'
$ws.Range("E50").Value = '  -2.28%  '
Set-TextValue "D51" '2.15'
$ws.Range("E51").Value = '  -1.39%  '
